$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns and fix "de/del/el/los" -> "De/Del/El/Los" capitalization
# in state/municipality names (data cleaning pass).
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B11").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Mazapa De Madero"
$ws.Range("A45").Value = "Ciudad De México"
$ws.Range("A56").Value = "Estado De México"
$ws.Range("B56").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B58").Value = "Almoloya De Juárez"
$ws.Range("B62").Value = "Naucalpan De Juárez"
$ws.Range("B64").Value = "Tlalnepantla De Baz"
$ws.Range("A68").Value = "Guanajuato"
$ws.Range("B68").Value = "Apaseo El Alto"
$ws.Range("B69").Value = "Apaseo El Grande"
$ws.Range("B71").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B75").Value = "Purísima Del Rincón"
$ws.Range("B76").Value = "San Luis De La Paz"
$ws.Range("B77").Value = "Silao De La Victoria"
$ws.Range("B78").Value = "Valle De Santiago"
$ws.Range("B81").Value = "Acapulco De Juárez"
$ws.Range("B83").Value = "Alcozauca De Guerero"
$ws.Range("B85").Value = "Atenango Del Río"
$ws.Range("B86").Value = "Atlamajalcingo Del Monte"
$ws.Range("B87").Value = "Ayutla De Los Libres"
$ws.Range("B88").Value = "Chilapa De Álvarez"
$ws.Range("B91").Value = "Iguala De La Independencia"
$ws.Range("B93").Value = "Mártir De Cuilapan"
$ws.Range("B99").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B108").Value = "Pachuca De Soto"
$ws.Range("B109").Value = "Tenango De Doria"
$ws.Range("B111").Value = "Autlán De Navarro"
$ws.Range("B118").Value = "Lagos De Moreno"
$ws.Range("B121").Value = "Ojuelos De Jalisco"
$ws.Range("B123").Value = "Tizapán El Alto"
$ws.Range("B124").Value = "Tlajomulco De Zúñiga"
$ws.Range("B127").Value = "Unión De Tula"
$ws.Range("B147").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B150").Value = "San Dionisio Del Mar"
$ws.Range("B154").Value = "San Pedro El Alto"
$ws.Range("B160").Value = "Tlalixtac De Cabrera"
$ws.Range("B161").Value = "Totontepec Villa De Morelos"
$ws.Range("B162").Value = "Villa Talea De Castro"
$ws.Range("B175").Value = "Izúcar De Matamoros"
$ws.Range("B183").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B202").Value = "Villa De Reyes"
$ws.Range("B227").Value = "Cazones De Herrera"
$ws.Range("B233").Value = "Ignacio De La Llave"
$ws.Range("B234").Value = "Ixhuatlán De Madero"
$ws.Range("B239").Value = "Martínez De La Torre"

# Drop the trailing footnote/source rows (259-263); data now ends at row 257.
$ws.Rows("259:263").Delete()

